$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,3,4) get cyclically shifted down by one:
# old row 2 -> new row 3
# old row 3 -> new row 4
# old row 4 -> new row 2
# Capture original values for the columns that change (D, J, K, L, M, O, P)

$origD2 = $ws.Range("D2").Value2
$origJ2 = $ws.Range("J2").Value2
$origK2 = $ws.Range("K2").Value2
$origL2 = $ws.Range("L2").Value2
$origM2 = $ws.Range("M2").Value2
$origO2 = $ws.Range("O2").Value2
$origP2 = $ws.Range("P2").Value2

$origD3 = $ws.Range("D3").Value2
$origJ3 = $ws.Range("J3").Value2
$origK3 = $ws.Range("K3").Value2
$origL3 = $ws.Range("L3").Value2
$origM3 = $ws.Range("M3").Value2
$origO3 = $ws.Range("O3").Value2
$origP3 = $ws.Range("P3").Value2

$origD4 = $ws.Range("D4").Value2
$origJ4 = $ws.Range("J4").Value2
$origK4 = $ws.Range("K4").Value2
$origL4 = $ws.Range("L4").Value2
$origM4 = $ws.Range("M4").Value2
$origO4 = $ws.Range("O4").Value2
$origP4 = $ws.Range("P4").Value2

# Row 2 gets old row 4's values
$ws.Range("D2").Value = $origD4
$ws.Range("J2").Value = $origJ4
$ws.Range("K2").Value = $origK4
$ws.Range("L2").Value = $origL4
$ws.Range("M2").Value = $origM4
$ws.Range("O2").Value = $origO4
$ws.Range("P2").Value = $origP4

# Row 3 gets old row 2's values
$ws.Range("D3").Value = $origD2
$ws.Range("J3").Value = $origJ2
$ws.Range("K3").Value = $origK2
$ws.Range("L3").Value = $origL2
$ws.Range("M3").Value = $origM2
$ws.Range("O3").Value = $origO2
$ws.Range("P3").Value = $origP2

# Row 4 gets old row 3's values
$ws.Range("D4").Value = $origD3
$ws.Range("J4").Value = $origJ3
$ws.Range("K4").Value = $origK3
$ws.Range("L4").Value = $origL3
$ws.Range("M4").Value = $origM3
$ws.Range("O4").Value = $origO3
$ws.Range("P4").Value = $origP3
